# The underlying source re-sorted/renumbered the observation rows for this
# species report. Rows 3,4,6,7,8,9,10,11,12 end up with each other's data
# (row 5 — and the header/row 2 — are untouched); column layout is unchanged.
# This is a pure permutation of whole rows, made of two disjoint cycles:
#   (3 12 7 8 10) and (4 11 9 6)
# i.e. new row 3 gets what used to be row 12's data, new row 12 gets what
# used to be row 7's data, and so on.
#
# Because it's a rotation (not independent edits), every source cell has to
# be read BEFORE any destination cell is written, otherwise an early write
# would clobber data a later step still needs to read. So: snapshot first,
# then write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> old row that its data comes from
$rowMap = @{
    3  = 12
    4  = 11
    6  = 4
    7  = 8
    8  = 10
    9  = 6
    10 = 3
    11 = 9
    12 = 7
}

# Only these columns ever differ between the affected rows; every other
# populated column (C, L, N, T, U, V, W, Y, AD, AE, AG, AT, AW, AX, AY, ...)
# holds the same value in all of them, so it is unaffected by the shuffle.
$cols = @("A","B","D","E","F","G","H","I","J","K","Q","R","S","Z","AB","AC")

# 1) Snapshot every affected cell's current value.
$snap = @{}
foreach ($r in $rowMap.Values) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        if (-not $snap.ContainsKey($addr)) {
            $snap[$addr] = $ws.Range($addr).Value2
        }
    }
}

# 2) Write the snapshotted values into their new row positions.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $snap["$col$oldRow"]
    }
}
